# Scheduled-runner update: refresh market-price-derived profit figures
# (currentAveragePrice*, Leve price/profit columns H:N) across all eight
# crafting-class sheets, matching the latest market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 226.26666
$ws.Range("I5").Value = 129.55556
$ws.Range("J5").Value = 371.33334
$ws.Range("K5").Value = 129.55556
$ws.Range("L5").Value = 371.33334
$ws.Range("M5").Value = -14.55556000000001
$ws.Range("N5").Value = -601.33334
$ws.Range("H9").Value = 6784.2666
$ws.Range("I9").Value = 9193.091
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 9193.091
$ws.Range("L9").Value = 160
$ws.Range("M9").Value = -9024.091
$ws.Range("N9").Value = -498
$ws.Range("H17").Value = 1842.5
$ws.Range("J17").Value = 1879.1305
$ws.Range("L17").Value = 5637.3915
$ws.Range("N17").Value = -5973.3915
$ws.Range("H39").Value = 1983.6
$ws.Range("I39").Value = 1048.375
$ws.Range("J39").Value = 5724.5
$ws.Range("K39").Value = 3145.125
$ws.Range("L39").Value = 17173.5
$ws.Range("M39").Value = -2849.125
$ws.Range("N39").Value = -17765.5
$ws.Range("H51").Value = 8627.177
$ws.Range("J51").Value = 6644.1333
$ws.Range("L51").Value = 6644.1333
$ws.Range("N51").Value = -7612.1333
$ws.Range("H70").Value = 6381.68
$ws.Range("I70").Value = 2415.7896
$ws.Range("K70").Value = 7247.3688
$ws.Range("M70").Value = -6977.3688
$ws.Range("H73").Value = 6381.68
$ws.Range("I73").Value = 2415.7896
$ws.Range("K73").Value = 7247.3688
$ws.Range("M73").Value = -6311.3688
$ws.Range("H74").Value = 5630.3335
$ws.Range("J74").Value = 5862.88
$ws.Range("L74").Value = 5862.88
$ws.Range("N74").Value = -7734.88
$ws.Range("H77").Value = 5630.3335
$ws.Range("J77").Value = 5862.88
$ws.Range("L77").Value = 29314.4
$ws.Range("N77").Value = -38674.4
$ws.Range("H97").Value = 4326
$ws.Range("J97").Value = 4326
$ws.Range("L97").Value = 12978
$ws.Range("N97").Value = -13970
$ws.Range("H100").Value = 18110.4
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H112").Value = 5296.607
$ws.Range("J112").Value = 5835
$ws.Range("L112").Value = 17505
$ws.Range("N112").Value = -19721
$ws.Range("H125").Value = 2447.6
$ws.Range("J125").Value = 2195.7693
$ws.Range("L125").Value = 19761.9237
$ws.Range("N125").Value = -24681.9237
$ws.Range("H129").Value = 1663.3846
$ws.Range("I129").Value = 1002.8889
$ws.Range("J129").Value = 3149.5
$ws.Range("K129").Value = 3008.6667
$ws.Range("L129").Value = 9448.5
$ws.Range("M129").Value = 1991.3333
$ws.Range("N129").Value = -19448.5
$ws.Range("H131").Value = 2772.6191
$ws.Range("I131").Value = 1412.5
$ws.Range("K131").Value = 4237.5
$ws.Range("M131").Value = 802.5
$ws.Range("H132").Value = 2806.4678
$ws.Range("I132").Value = 2656.3728
$ws.Range("K132").Value = 7969.1184
$ws.Range("M132").Value = -5439.1184
$ws.Range("H138").Value = 5281.9355
$ws.Range("I138").Value = 5579.857
$ws.Range("K138").Value = 16739.571
$ws.Range("M138").Value = -11599.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1566.5333
$ws.Range("I2").Value = 1140.4
$ws.Range("K2").Value = 1140.4
$ws.Range("M2").Value = -1027.4
$ws.Range("H32").Value = 2677.9092
$ws.Range("I32").Value = 1289.8572
$ws.Range("K32").Value = 1289.8572
$ws.Range("M32").Value = -1002.8572
$ws.Range("H61").Value = 4135.4116
$ws.Range("I61").Value = 3923.2307
$ws.Range("J61").Value = 4825
$ws.Range("K61").Value = 3923.2307
$ws.Range("L61").Value = 4825
$ws.Range("M61").Value = -3711.2307
$ws.Range("N61").Value = -5249
$ws.Range("H116").Value = 1566.5333
$ws.Range("I116").Value = 1140.4
$ws.Range("K116").Value = 1140.4
$ws.Range("M116").Value = 1153.6
$ws.Range("H132").Value = 2092.6667
$ws.Range("I132").Value = 1820.5883
$ws.Range("K132").Value = 5461.7649
$ws.Range("M132").Value = -2931.7649
$ws.Range("H136").Value = 4135.4116
$ws.Range("I136").Value = 3923.2307
$ws.Range("J136").Value = 4825
$ws.Range("K136").Value = 11769.6921
$ws.Range("L136").Value = 14475
$ws.Range("M136").Value = -9219.6921
$ws.Range("N136").Value = -19575

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1566.5333
$ws.Range("I3").Value = 1140.4
$ws.Range("K3").Value = 1140.4
$ws.Range("M3").Value = -1026.4
$ws.Range("H99").Value = 2360.375
$ws.Range("I99").Value = 783.8
$ws.Range("K99").Value = 783.8
$ws.Range("M99").Value = 714.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6608.8184
$ws.Range("I86").Value = 6083.222
$ws.Range("K86").Value = 6083.222
$ws.Range("M86").Value = -4960.222
$ws.Range("H89").Value = 6608.8184
$ws.Range("I89").Value = 6083.222
$ws.Range("K89").Value = 30416.11
$ws.Range("M89").Value = -24800.11
$ws.Range("H122").Value = 4133.1665
$ws.Range("I122").Value = 2726.125
$ws.Range("K122").Value = 8178.375
$ws.Range("M122").Value = -5728.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1381.875
$ws.Range("I122").Value = 263.75
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 2373.75
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = 76.25
$ws.Range("N122").Value = -27400
$ws.Range("H129").Value = 12502470
$ws.Range("I129").Value = 27778540
$ws.Range("J129").Value = 9806692
$ws.Range("K129").Value = 83335620
$ws.Range("L129").Value = 29420076
$ws.Range("M129").Value = -83330620
$ws.Range("N129").Value = -29430076
$ws.Range("H131").Value = 31748054
$ws.Range("I131").Value = 27779628
$ws.Range("J131").Value = 33335424
$ws.Range("K131").Value = 83338884
$ws.Range("L131").Value = 100006272
$ws.Range("M131").Value = -83333844
$ws.Range("N131").Value = -100016352
$ws.Range("H134").Value = 13390.392
$ws.Range("I134").Value = 1333.1666
$ws.Range("K134").Value = 3999.4998
$ws.Range("M134").Value = 1070.5002
$ws.Range("H139").Value = 12831429
$ws.Range("I139").Value = 23818798
$ws.Range("J139").Value = 12833.333
$ws.Range("K139").Value = 71456394
$ws.Range("L139").Value = 38499.999
$ws.Range("M139").Value = -71451254
$ws.Range("N139").Value = -48779.999
$ws.Range("H140").Value = 31261272
$ws.Range("I140").Value = 62505292
$ws.Range("K140").Value = 187515876
$ws.Range("M140").Value = -187510696

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 36474.5
$ws.Range("J26").Value = 35632.668
$ws.Range("L26").Value = 35632.668
$ws.Range("N26").Value = -36192.668
$ws.Range("H50").Value = 36474.5
$ws.Range("J50").Value = 35632.668
$ws.Range("L50").Value = 35632.668
$ws.Range("N50").Value = -36628.668
$ws.Range("H97").Value = 2773.6
$ws.Range("I97").Value = 561.6
$ws.Range("K97").Value = 561.6
$ws.Range("M97").Value = -65.60000000000002
$ws.Range("H102").Value = 9238.546
$ws.Range("I102").Value = 10577.846
$ws.Range("J102").Value = 4264
$ws.Range("K102").Value = 10577.846
$ws.Range("L102").Value = 4264
$ws.Range("M102").Value = -8955.846
$ws.Range("N102").Value = -7508
$ws.Range("H122").Value = 3961.9062
$ws.Range("I122").Value = 3420.8823
$ws.Range("K122").Value = 10262.6469
$ws.Range("M122").Value = -7812.6469
$ws.Range("H132").Value = 2712
$ws.Range("I132").Value = 2055.5
$ws.Range("J132").Value = 3499.8
$ws.Range("K132").Value = 6166.5
$ws.Range("L132").Value = 10499.4
$ws.Range("M132").Value = -3636.5
$ws.Range("N132").Value = -15559.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1629.421
$ws.Range("J46").Value = 1980.5834
$ws.Range("L46").Value = 1980.5834
$ws.Range("N46").Value = -2356.5834
$ws.Range("H55").Value = 315.23077
$ws.Range("J55").Value = 397.75
$ws.Range("L55").Value = 397.75
$ws.Range("N55").Value = -743.75
$ws.Range("H136").Value = 2212.7368
$ws.Range("I136").Value = 2042.4546
$ws.Range("K136").Value = 6127.3638
$ws.Range("M136").Value = -3577.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9972.166999999999
$ws.Range("I81").Value = 9972.166999999999
$ws.Range("K81").Value = 19944.334
$ws.Range("M81").Value = -18883.334
$ws.Range("H84").Value = 9972.166999999999
$ws.Range("I84").Value = 9972.166999999999
$ws.Range("K84").Value = 99721.67
$ws.Range("M84").Value = -94417.67
$ws.Range("H126").Value = 2586.875
$ws.Range("I126").Value = 2556.4285
$ws.Range("K126").Value = 7669.2855
$ws.Range("M126").Value = -5199.2855
